$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 352 (shifts existing rows 352-434 down to 354-436)
$ws.Rows.Item(352).Insert()
$ws.Rows.Item(352).Insert()

# New row 352: Albahaca "Primera" quality entry for date 44641 (2022-03-21)
$ws.Range("A352").Value = 6
$ws.Range("B352").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C352").Value = "Metropolitana"
$ws.Range("D352").Value = 44641
$ws.Range("E352").Value = 13
$ws.Range("F352").Value = 100112052
$ws.Range("G352").Value = "Albahaca"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 190
$ws.Range("K352").Value = 2500
$ws.Range("L352").Value = 2500
$ws.Range("M352").Value = 2500
$ws.Range("N352").Value = "`$/docena de matas"
$ws.Range("O352").Value = "Región Metropolitana"
$ws.Range("P352").Value = 417
$ws.Range("Q352").Value = 6
$ws.Range("R352").Value = "Hortaliza"

# New row 353: Albahaca "Segunda" quality entry for date 44641 (2022-03-21)
$ws.Range("A353").Value = 6
$ws.Range("B353").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C353").Value = "Metropolitana"
$ws.Range("D353").Value = 44641
$ws.Range("E353").Value = 13
$ws.Range("F353").Value = 100112052
$ws.Range("G353").Value = "Albahaca"
$ws.Range("H353").Value = "Sin especificar"
$ws.Range("I353").Value = "Segunda"
$ws.Range("J353").Value = 130
$ws.Range("K353").Value = 2000
$ws.Range("L353").Value = 2000
$ws.Range("M353").Value = 2000
$ws.Range("N353").Value = "`$/docena de matas"
$ws.Range("O353").Value = "Región Metropolitana"
$ws.Range("P353").Value = 333
$ws.Range("Q353").Value = 6
$ws.Range("R353").Value = "Hortaliza"

# Apply the date number format used by the other Fecha (column D) cells
$ws.Range("D352").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D353").NumberFormat = "YYYY-MM-DD HH:MM:SS"
